$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 123.9889473333333
$ws.Cells.Item(2, 8).Value = 371.966842
$ws.Cells.Item(2, 9).Value = 0.1954468191201633
$ws.Cells.Item(2, 10).Value = 0.2134822208566356
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.114644
$ws.Cells.Item(2, 14).Value = 0.343932
$ws.Cells.Item(2, 15).Value = 0.107929744556041
$ws.Cells.Item(2, 16).Value = 0.1536053250846448
$ws.Cells.Item(2, 17).Value = 14.21458887808267
$ws.Cells.Item(2, 18).Value = 127.931299902744
$ws.Cells.Item(2, 19).Value = 0.02109452526192997
$ws.Cells.Item(2, 20).Value = 0.03279200593447545

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 123.9889473333333
$ws.Cells.Item(3, 8).Value = 371.966842
$ws.Cells.Item(3, 9).Value = 0.1954468191201633
$ws.Cells.Item(3, 10).Value = 0.2134822208566356
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 13).Value = 0.9475654999999999
$ws.Cells.Item(3, 14).Value = 1.895131
$ws.Cells.Item(3, 15).Value = 0.8920702554439589
$ws.Cells.Item(3, 16).Value = 0.8463946749153552
$ws.Cells.Item(3, 17).Value = 117.4876488743837
$ws.Cells.Item(3, 18).Value = 704.9258932463021
$ws.Cells.Item(3, 19).Value = 0.1743522938582333
$ws.Cells.Item(3, 20).Value = 0.1806902149221602

$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 211.2281266666666
$ws.Cells.Item(4, 8).Value = 633.6843799999999
$ws.Cells.Item(4, 9).Value = 0.3329640774731549
$ws.Cells.Item(4, 10).Value = 0.3636892687455195
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.114644
$ws.Cells.Item(4, 14).Value = 0.343932
$ws.Cells.Item(4, 15).Value = 0.107929744556041
$ws.Cells.Item(4, 16).Value = 0.1536053250846448
$ws.Cells.Item(4, 17).Value = 24.21603735357333
$ws.Cells.Item(4, 18).Value = 217.94433618216
$ws.Cells.Item(4, 19).Value = 0.03593672782801545
$ws.Cells.Item(4, 20).Value = 0.05586460835545227

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 211.2281266666666
$ws.Cells.Item(5, 8).Value = 633.6843799999999
$ws.Cells.Item(5, 9).Value = 0.3329640774731549
$ws.Cells.Item(5, 10).Value = 0.3636892687455195
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 13).Value = 0.9475654999999999
$ws.Cells.Item(5, 14).Value = 1.895131
$ws.Cells.Item(5, 15).Value = 0.8920702554439589
$ws.Cells.Item(5, 16).Value = 0.8463946749153552
$ws.Cells.Item(5, 17).Value = 200.1524854589633
$ws.Cells.Item(5, 18).Value = 1200.91491275378
$ws.Cells.Item(5, 19).Value = 0.2970273496451394
$ws.Cells.Item(5, 20).Value = 0.3078246603900672

$ws.Cells.Item(6, 1).Value = "M1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 77.63463966666666
$ws.Cells.Item(6, 8).Value = 232.903919
$ws.Cells.Item(6, 9).Value = 0.1223773868778609
$ws.Cells.Item(6, 10).Value = 0.1336701024397599
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.114644
$ws.Cells.Item(6, 14).Value = 0.343932
$ws.Cells.Item(6, 15).Value = 0.107929744556041
$ws.Cells.Item(6, 16).Value = 0.1536053250846448
$ws.Cells.Item(6, 17).Value = 8.900345629945333
$ws.Cells.Item(6, 18).Value = 80.103110669508
$ws.Cells.Item(6, 19).Value = 0.01320816010516332
$ws.Cells.Item(6, 20).Value = 0.02053243953935709

$ws.Cells.Item(7, 1).Value = "M1"
$ws.Cells.Item(7, 2).Value = "Gnas"
$ws.Cells.Item(7, 3).Value = "Lhcgr"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 77.63463966666666
$ws.Cells.Item(7, 8).Value = 232.903919
$ws.Cells.Item(7, 9).Value = 0.1223773868778609
$ws.Cells.Item(7, 10).Value = 0.1336701024397599
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.9475654999999999
$ws.Cells.Item(7, 14).Value = 1.895131
$ws.Cells.Item(7, 15).Value = 0.8920702554439589
$ws.Cells.Item(7, 16).Value = 0.8463946749153552
$ws.Cells.Item(7, 17).Value = 73.56390615306482
$ws.Cells.Item(7, 18).Value = 441.383436918389
$ws.Cells.Item(7, 19).Value = 0.1091692267726975
$ws.Cells.Item(7, 20).Value = 0.1131376629004028

$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Gnas"
$ws.Cells.Item(8, 3).Value = "Lhcgr"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 60.75256733333333
$ws.Cells.Item(8, 8).Value = 182.257702
$ws.Cells.Item(8, 9).Value = 0.09576576214298858
$ws.Cells.Item(8, 10).Value = 0.1046028156219013
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.114644
$ws.Cells.Item(8, 14).Value = 0.343932
$ws.Cells.Item(8, 15).Value = 0.107929744556041
$ws.Cells.Item(8, 16).Value = 0.1536053250846448
$ws.Cells.Item(8, 17).Value = 6.964917329362667
$ws.Cells.Item(8, 18).Value = 62.684255964264
$ws.Cells.Item(8, 19).Value = 0.01033597424530734
$ws.Cells.Item(8, 20).Value = 0.0160675494983713

$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Gnas"
$ws.Cells.Item(9, 3).Value = "Lhcgr"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 60.75256733333333
$ws.Cells.Item(9, 8).Value = 182.257702
$ws.Cells.Item(9, 9).Value = 0.09576576214298858
$ws.Cells.Item(9, 10).Value = 0.1046028156219013
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.9475654999999999
$ws.Cells.Item(9, 14).Value = 1.895131
$ws.Cells.Item(9, 15).Value = 0.8920702554439589
$ws.Cells.Item(9, 16).Value = 0.8463946749153552
$ws.Cells.Item(9, 17).Value = 57.56703684149366
$ws.Cells.Item(9, 18).Value = 345.402221048962
$ws.Cells.Item(9, 19).Value = 0.08542978789768124
$ws.Cells.Item(9, 20).Value = 0.08853526612352995

$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Gnas"
$ws.Cells.Item(10, 3).Value = "Lhcgr"
$ws.Cells.Item(10, 4).Value = "FAPs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 160.7828525
$ws.Cells.Item(10, 8).Value = 321.565705
$ws.Cells.Item(10, 9).Value = 0.2534459543858325
$ws.Cells.Item(10, 10).Value = 0.1845555923361839
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.114644
$ws.Cells.Item(10, 14).Value = 0.343932
$ws.Cells.Item(10, 15).Value = 0.107929744556041
$ws.Cells.Item(10, 16).Value = 0.1536053250846448
$ws.Cells.Item(10, 17).Value = 18.43278934201
$ws.Cells.Item(10, 18).Value = 110.59673605206
$ws.Cells.Item(10, 19).Value = 0.02735435711562491
$ws.Cells.Item(10, 20).Value = 0.02834872175698871

$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Gnas"
$ws.Cells.Item(11, 3).Value = "Lhcgr"
$ws.Cells.Item(11, 4).Value = "sCs"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 160.7828525
$ws.Cells.Item(11, 8).Value = 321.565705
$ws.Cells.Item(11, 9).Value = 0.2534459543858325
$ws.Cells.Item(11, 10).Value = 0.1845555923361839
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.9475654999999999
$ws.Cells.Item(11, 14).Value = 1.895131
$ws.Cells.Item(11, 15).Value = 0.8920702554439589
$ws.Cells.Item(11, 16).Value = 0.8463946749153552
$ws.Cells.Item(11, 17).Value = 152.3522840205887
$ws.Cells.Item(11, 18).Value = 609.409136082355
$ws.Cells.Item(11, 19).Value = 0.2260915972702075
$ws.Cells.Item(11, 20).Value = 0.1562068705791952
